$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Cells.Item(2, 1).Value = "100%>gen/cap>95%"
$ws.Cells.Item(2, 2).Value = [double]"0"
$ws.Cells.Item(2, 8).Value = [double]"0.07583158525043154"

$ws.Cells.Item(3, 1).Value = "105%>gen/cap>100%"
$ws.Cells.Item(3, 2).Value = [double]"0.003147905835586739"
$ws.Cells.Item(3, 3).Value = [double]"0.0007066694810891484"
$ws.Cells.Item(3, 4).Value = [double]"4.539406795542709"
$ws.Cells.Item(3, 5).Value = [double]"2.175451737738439e-05"
$ws.Cells.Item(3, 6).Value = [double]"0.001762855890208767"
$ws.Cells.Item(3, 7).Value = [double]"0.00453295578096471"
$ws.Cells.Item(3, 8).Value = [double]"0.07897949108601829"

$ws.Cells.Item(4, 1).Value = "110%>gen/cap>105%"
$ws.Cells.Item(4, 2).Value = [double]"0.002122880227381964"
$ws.Cells.Item(4, 3).Value = [double]"0.001195267588934897"
$ws.Cells.Item(4, 4).Value = [double]"0.5719720849629968"
$ws.Cells.Item(4, 5).Value = [double]"0.0007576066091514218"
$ws.Cells.Item(4, 6).Value = [double]"-0.0002198068331577922"
$ws.Cells.Item(4, 7).Value = [double]"0.00446556728792172"
$ws.Cells.Item(4, 8).Value = [double]"0.07795446547781351"

$ws.Cells.Item(5, 1).Value = "115%>gen/cap>110%"
$ws.Cells.Item(5, 2).Value = [double]"0.006818716880059835"
$ws.Cells.Item(5, 3).Value = [double]"0.00168137056793241"
$ws.Cells.Item(5, 4).Value = [double]"6.398669531437461"
$ws.Cells.Item(5, 5).Value = [double]"0.0003222463161229662"
$ws.Cells.Item(5, 6).Value = [double]"0.003523284384575951"
$ws.Cells.Item(5, 7).Value = [double]"0.01011414937554372"
$ws.Cells.Item(5, 8).Value = [double]"0.08265030213049138"

$ws.Cells.Item(6, 1).Value = "120%>gen/cap>115%"
$ws.Cells.Item(6, 2).Value = [double]"0.01657114076009747"
$ws.Cells.Item(6, 3).Value = [double]"0.01294462043818133"
$ws.Cells.Item(6, 4).Value = [double]"3.950586022558934"
$ws.Cells.Item(6, 5).Value = [double]"0.259782470326692"
$ws.Cells.Item(6, 6).Value = [double]"-0.00879989537445414"
$ws.Cells.Item(6, 7).Value = [double]"0.04194217689464907"
$ws.Cells.Item(6, 8).Value = [double]"0.09240272601052901"

$ws.Cells.Item(7, 1).Value = "125%>gen/cap>120%"
$ws.Cells.Item(7, 2).Value = [double]"0.1257131355860988"
$ws.Cells.Item(7, 3).Value = [double]"0.07310079967184438"
$ws.Cells.Item(7, 4).Value = [double]"9.136965276924572"
$ws.Cells.Item(7, 5).Value = [double]"9.080892016264132"
$ws.Cells.Item(7, 6).Value = [double]"0.05607326066043971"
$ws.Cells.Item(7, 7).Value = [double]"0.07310079967184438"
$ws.Cells.Item(7, 8).Value = [double]"0.2015447208365304"

$ws.Cells.Item(8, 1).Value = "130%>gen/cap>125%"
$ws.Cells.Item(8, 2).Value = [double]"0.1316002697888477"
$ws.Cells.Item(8, 3).Value = [double]"0.08280597133041581"
$ws.Cells.Item(8, 4).Value = [double]"5.780996977241826"
$ws.Cells.Item(8, 5).Value = [double]"5.726457404677218"
$ws.Cells.Item(8, 6).Value = [double]"0.05453957256460788"
$ws.Cells.Item(8, 7).Value = [double]"0.08280597133041581"
$ws.Cells.Item(8, 8).Value = [double]"0.2074318550392792"

$ws.Cells.Item(9, 1).Value = "135%>gen/cap>130%"
$ws.Cells.Item(9, 2).Value = [double]"0.1458089334744864"
$ws.Cells.Item(9, 3).Value = [double]"0.08837208485250299"
$ws.Cells.Item(9, 4).Value = [double]"7.038364858876592"
$ws.Cells.Item(9, 5).Value = [double]"6.975702409845225"
$ws.Cells.Item(9, 6).Value = [double]"0.0626624490313665"
$ws.Cells.Item(9, 7).Value = [double]"0.08837208485250299"
$ws.Cells.Item(9, 8).Value = [double]"0.221640518724918"

$ws.Cells.Item(10, 1).Value = "20%>gen/cap"
$ws.Cells.Item(10, 2).Value = [double]"-0.07583158525043154"
$ws.Cells.Item(10, 3).Value = [double]"0.0005217811936168974"
$ws.Cells.Item(10, 4).Value = [double]"-158.9312206264254"
$ws.Cells.Item(10, 5).Value = [double]"0"
$ws.Cells.Item(10, 6).Value = [double]"-0.07685425999964621"
$ws.Cells.Item(10, 7).Value = [double]"-0.0748089105012169"
$ws.Cells.Item(10, 8).Value = [double]"0"

$ws.Cells.Item(11, 1).Value = "25%>gen/cap>20%"
$ws.Cells.Item(11, 2).Value = [double]"-0.03208944258288336"
$ws.Cells.Item(11, 3).Value = [double]"0.0005533304317203182"
$ws.Cells.Item(11, 4).Value = [double]"-61.45664515634908"
$ws.Cells.Item(11, 5).Value = [double]"0"
$ws.Cells.Item(11, 6).Value = [double]"-0.03317395284494203"
$ws.Cells.Item(11, 7).Value = [double]"-0.0310049323208247"
$ws.Cells.Item(11, 8).Value = [double]"0.04374214266754818"

$ws.Cells.Item(12, 1).Value = "30%>gen/cap>25%"
$ws.Cells.Item(12, 2).Value = [double]"-0.02191124259456369"
$ws.Cells.Item(12, 3).Value = [double]"0.0005540901271342327"
$ws.Cells.Item(12, 4).Value = [double]"-41.59653024125657"
$ws.Cells.Item(12, 5).Value = [double]"8.98036756130432e-231"
$ws.Cells.Item(12, 6).Value = [double]"-0.02299724183174583"
$ws.Cells.Item(12, 7).Value = [double]"-0.02082524335738156"
$ws.Cells.Item(12, 8).Value = [double]"0.05392034265586785"

$ws.Cells.Item(13, 1).Value = "35%>gen/cap>30%"
$ws.Cells.Item(13, 2).Value = [double]"-0.01885740209498997"
$ws.Cells.Item(13, 3).Value = [double]"0.0005371419629511076"
$ws.Cells.Item(13, 4).Value = [double]"-36.61414020383948"
$ws.Cells.Item(13, 5).Value = [double]"4.607702551601126e-199"
$ws.Cells.Item(13, 6).Value = [double]"-0.01991018346787593"
$ws.Cells.Item(13, 7).Value = [double]"-0.01780462072210401"
$ws.Cells.Item(13, 8).Value = [double]"0.05697418315544157"

$ws.Cells.Item(14, 1).Value = "40%>gen/cap>35%"
$ws.Cells.Item(14, 2).Value = [double]"-0.01571388199847679"
$ws.Cells.Item(14, 3).Value = [double]"0.0005362836766911769"
$ws.Cells.Item(14, 4).Value = [double]"-30.40523081510505"
$ws.Cells.Item(14, 5).Value = [double]"4.750876254368663e-144"
$ws.Cells.Item(14, 6).Value = [double]"-0.01676498115428728"
$ws.Cells.Item(14, 7).Value = [double]"-0.0146627828426663"
$ws.Cells.Item(14, 8).Value = [double]"0.06011770325195476"

$ws.Cells.Item(15, 1).Value = "45%>gen/cap>40%"
$ws.Cells.Item(15, 2).Value = [double]"-0.009437293247196627"
$ws.Cells.Item(15, 3).Value = [double]"0.0005339264070832462"
$ws.Cells.Item(15, 4).Value = [double]"-18.70090529324256"
$ws.Cells.Item(15, 5).Value = [double]"5.204580603409424e-46"
$ws.Cells.Item(15, 6).Value = [double]"-0.01048377222912063"
$ws.Cells.Item(15, 7).Value = [double]"-0.00839081426527262"
$ws.Cells.Item(15, 8).Value = [double]"0.06639429200323492"

$ws.Cells.Item(16, 1).Value = "50%>gen/cap>45%"
$ws.Cells.Item(16, 2).Value = [double]"-0.007473542049845478"
$ws.Cells.Item(16, 3).Value = [double]"0.0005296077436073473"
$ws.Cells.Item(16, 4).Value = [double]"-15.16657722950536"
$ws.Cells.Item(16, 5).Value = [double]"9.698244986982895e-28"
$ws.Cells.Item(16, 6).Value = [double]"-0.008511556598480925"
$ws.Cells.Item(16, 7).Value = [double]"-0.00643552750121003"
$ws.Cells.Item(16, 8).Value = [double]"0.06835804320058607"

$ws.Cells.Item(17, 1).Value = "55%>gen/cap>50%"
$ws.Cells.Item(17, 2).Value = [double]"-0.008369918839504498"
$ws.Cells.Item(17, 3).Value = [double]"0.0005516555121101367"
$ws.Cells.Item(17, 4).Value = [double]"-16.07704048253514"
$ws.Cells.Item(17, 5).Value = [double]"6.64934180852886e-35"
$ws.Cells.Item(17, 6).Value = [double]"-0.00945114631728048"
$ws.Cells.Item(17, 7).Value = [double]"-0.007288691361728513"
$ws.Cells.Item(17, 8).Value = [double]"0.06746166641092705"

$ws.Cells.Item(18, 1).Value = "60%>gen/cap>55%"
$ws.Cells.Item(18, 2).Value = [double]"-0.006834642358546209"
$ws.Cells.Item(18, 3).Value = [double]"0.0005560578762806755"
$ws.Cells.Item(18, 4).Value = [double]"-13.23104090931768"
$ws.Cells.Item(18, 5).Value = [double]"7.82188352252107e-21"
$ws.Cells.Item(18, 6).Value = [double]"-0.00792449833237065"
$ws.Cells.Item(18, 7).Value = [double]"-0.005744786384721768"
$ws.Cells.Item(18, 8).Value = [double]"0.06899694289188534"

$ws.Cells.Item(19, 1).Value = "65%>gen/cap>60%"
$ws.Cells.Item(19, 2).Value = [double]"-0.005569882318398344"
$ws.Cells.Item(19, 3).Value = [double]"0.0005452686004827275"
$ws.Cells.Item(19, 4).Value = [double]"-10.98765540548857"
$ws.Cells.Item(19, 5).Value = [double]"6.370248240575907e-16"
$ws.Cells.Item(19, 6).Value = [double]"-0.006638591664270381"
$ws.Cells.Item(19, 7).Value = [double]"-0.004501172972526309"
$ws.Cells.Item(19, 8).Value = [double]"0.0702617029320332"

$ws.Cells.Item(20, 1).Value = "70%>gen/cap>65%"
$ws.Cells.Item(20, 2).Value = [double]"-0.006432836614972464"
$ws.Cells.Item(20, 3).Value = [double]"0.0005562938338018325"
$ws.Cells.Item(20, 4).Value = [double]"-12.47649270275535"
$ws.Cells.Item(20, 5).Value = [double]"2.426023993181101e-17"
$ws.Cells.Item(20, 6).Value = [double]"-0.007523155047592909"
$ws.Cells.Item(20, 7).Value = [double]"-0.00534251818235202"
$ws.Cells.Item(20, 8).Value = [double]"0.06939874863545908"

$ws.Cells.Item(21, 1).Value = "75%>gen/cap>70%"
$ws.Cells.Item(21, 2).Value = [double]"-0.003788248138857311"
$ws.Cells.Item(21, 3).Value = [double]"0.0005662406734161853"
$ws.Cells.Item(21, 4).Value = [double]"-7.579529083267804"
$ws.Cells.Item(21, 5).Value = [double]"0.0001982447943883382"
$ws.Cells.Item(21, 6).Value = [double]"-0.004898062062330771"
$ws.Cells.Item(21, 7).Value = [double]"-0.002678434215383851"
$ws.Cells.Item(21, 8).Value = [double]"0.07204333711157423"

$ws.Cells.Item(22, 1).Value = "80%>gen/cap>75%"
$ws.Cells.Item(22, 2).Value = [double]"-0.001977077191380415"
$ws.Cells.Item(22, 3).Value = [double]"0.0005585592978857748"
$ws.Cells.Item(22, 4).Value = [double]"-4.489024130927394"
$ws.Cells.Item(22, 5).Value = [double]"0.1846580007453567"
$ws.Cells.Item(22, 6).Value = [double]"-0.003071835829200818"
$ws.Cells.Item(22, 7).Value = [double]"-0.0008823185535600123"
$ws.Cells.Item(22, 8).Value = [double]"0.07385450805905112"

$ws.Cells.Item(23, 1).Value = "85%>gen/cap>80%"
$ws.Cells.Item(23, 2).Value = [double]"-0.0007993037230699283"
$ws.Cells.Item(23, 3).Value = [double]"0.0005549203389480512"
$ws.Cells.Item(23, 4).Value = [double]"-1.913292878071929"
$ws.Cells.Item(23, 5).Value = [double]"0.2907939118841092"
$ws.Cells.Item(23, 6).Value = [double]"-0.001886930126893345"
$ws.Cells.Item(23, 7).Value = [double]"0.0002883226807534879"
$ws.Cells.Item(23, 8).Value = [double]"0.07503228152736162"

$ws.Cells.Item(24, 1).Value = "90%>gen/cap>85%"
$ws.Cells.Item(24, 2).Value = [double]"-0.001267578312711584"
$ws.Cells.Item(24, 3).Value = [double]"0.0005591556371085686"
$ws.Cells.Item(24, 4).Value = [double]"-2.209690571768653"
$ws.Cells.Item(24, 5).Value = [double]"0.03132938893094828"
$ws.Cells.Item(24, 6).Value = [double]"-0.002363505805619614"
$ws.Cells.Item(24, 7).Value = [double]"-0.0001716508198035535"
$ws.Cells.Item(24, 8).Value = [double]"0.07456400693771996"

$ws.Cells.Item(25, 1).Value = "95%>gen/cap>90%"
$ws.Cells.Item(25, 2).Value = [double]"-0.001588739919502577"
$ws.Cells.Item(25, 3).Value = [double]"0.0005419654112399709"
$ws.Cells.Item(25, 4).Value = [double]"-2.18971004877795"
$ws.Cells.Item(25, 5).Value = [double]"0.2856785743342236"
$ws.Cells.Item(25, 6).Value = [double]"-0.00265097512546101"
$ws.Cells.Item(25, 7).Value = [double]"-0.000526504713544143"
$ws.Cells.Item(25, 8).Value = [double]"0.07424284533092897"

$ws.Cells.Item(26, 1).Value = "gen/cap>135%"
$ws.Cells.Item(26, 2).Value = [double]"0.1441788494195033"
$ws.Cells.Item(26, 3).Value = [double]"0.0773544133941799"
$ws.Cells.Item(26, 4).Value = [double]"20.33289195488689"
$ws.Cells.Item(26, 5).Value = [double]"20.26428878817274"
$ws.Cells.Item(26, 6).Value = [double]"0.06860316671415534"
$ws.Cells.Item(26, 7).Value = [double]"0.0773544133941799"
$ws.Cells.Item(26, 8).Value = [double]"0.2200104346699348"
